$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 84 should visually/structurally match the existing data rows
# (79-83), which use style index 4 for columns A-E and style index 5
# (date format) for column F, plus a custom row height of 15.35.
# Copy formatting from the last existing data row (83) down to the new row 84.
$ws.Range("A83:E83").Copy()
$ws.Range("A84:E84").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("F83").Copy()
$ws.Range("F84").PasteSpecial(-4122)      # xlPasteFormats

$ws.Rows.Item(84).RowHeight = 15.35

# Populate the new customer record:
# name, tier, city, state_prov, country, created_date
$ws.Cells.Item(84, 1).Value = "Anthony Pigatto"
$ws.Cells.Item(84, 2).Value = "Wholesale"
$ws.Cells.Item(84, 3).Value = "La Grange"
$ws.Cells.Item(84, 4).Value = "IL"
$ws.Cells.Item(84, 5).Value = "USA"
$ws.Cells.Item(84, 6).Value = 45643

Write-Host "Added row 84:" $ws.Cells.Item(84, 1).Value()
